$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ENBW charging session): add Preis/kWh, fix Kennzeichen typo, normalize Kosten format
# "0.74" parses as a number in general format, so force text formatting first
# so it is stored the same way as the other text-like values in the sheet.
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "0.74"
$ws.Range("I2").Value = "PBSN23E"
$ws.Range("K2").Value = "58.91 EUR"

# Row 3 (Aral, Super): add missing Liter value, normalize Kosten format (comma -> dot)
$ws.Range("D3").Value = "42,63"
$ws.Range("K3").Value = "73.71 EUR"

# Row 4 (HEM, Super): normalize Kosten format (comma -> dot)
$ws.Range("K4").Value = "70.55 EUR"

# Row 5 (Aral, Super E10): add missing Liter value, normalize Kosten format (comma -> dot)
$ws.Range("D5").Value = "43,84"
$ws.Range("K5").Value = "73.17 EUR"

# Row 6 (Total, Waschkarte): normalize Kosten format (comma -> dot)
$ws.Range("K6").Value = "17.50 EUR"

# Row 7 (HEM, Super): normalize Kosten format (comma -> dot)
$ws.Range("K7").Value = "61.75 EUR"

# Row 8 (Total, Super E10): normalize Kosten format (comma -> dot)
$ws.Range("K8").Value = "73.32 EUR"
